$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1510
$ws.Range("I6").Value = 190
$ws.Range("J6").Value = 2500
$ws.Range("K6").Value = 570
$ws.Range("L6").Value = 7500
$ws.Range("M6").Value = -458
$ws.Range("N6").Value = -7724
$ws.Range("H13").Value = 8646.857
$ws.Range("I13").Value = 2900
$ws.Range("J13").Value = 17528.363
$ws.Range("K13").Value = 2900
$ws.Range("L13").Value = 17528.363
$ws.Range("M13").Value = -2731
$ws.Range("N13").Value = -17866.363
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("H39").Value = 450.2857
$ws.Range("J39").Value = 1009.8889
$ws.Range("L39").Value = 3029.6667
$ws.Range("N39").Value = -3621.6667
$ws.Range("H76").Value = 2971.75
$ws.Range("I76").Value = 2795.5715
$ws.Range("J76").Value = 4205
$ws.Range("K76").Value = 2795.5715
$ws.Range("L76").Value = 4205
$ws.Range("M76").Value = -2480.5715
$ws.Range("N76").Value = -4835
$ws.Range("H79").Value = 2971.75
$ws.Range("I79").Value = 2795.5715
$ws.Range("J79").Value = 4205
$ws.Range("K79").Value = 2795.5715
$ws.Range("L79").Value = 4205
$ws.Range("M79").Value = -1703.5715
$ws.Range("N79").Value = -6389
$ws.Range("H132").Value = 3510775.8
$ws.Range("I132").Value = 4083651
$ws.Range("J132").Value = 1914.875
$ws.Range("K132").Value = 12250953
$ws.Range("L132").Value = 5744.625
$ws.Range("M132").Value = -12248423
$ws.Range("N132").Value = -10804.625

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 73334.336
$ws.Range("I11").Value = 65000
$ws.Range("J11").Value = 77501.5
$ws.Range("K11").Value = 65000
$ws.Range("L11").Value = 77501.5
$ws.Range("M11").Value = -64856
$ws.Range("N11").Value = -77789.5
$ws.Range("H28").Value = 4840
$ws.Range("I28").Value = 4840
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 4840
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -4648
$ws.Range("N28").ClearContents()
$ws.Range("H70").Value = 38750
$ws.Range("J70").Value = 38750
$ws.Range("L70").Value = 38750
$ws.Range("N70").Value = -39290
$ws.Range("H73").Value = 38750
$ws.Range("J73").Value = 38750
$ws.Range("L73").Value = 38750
$ws.Range("N73").Value = -40622
$ws.Range("H99").Value = 4840
$ws.Range("I99").Value = 4840
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 4840
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1845
$ws.Range("N99").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 36502.332
$ws.Range("J15").Value = 70004.664
$ws.Range("L15").Value = 70004.664
$ws.Range("N15").Value = -70458.664
$ws.Range("H128").Value = 2000
$ws.Range("I128").Value = 2000
$ws.Range("K128").Value = 6000
$ws.Range("M128").Value = -3510
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H54").Value = 44051
$ws.Range("J54").Value = 44051
$ws.Range("L54").Value = 44051
$ws.Range("N54").Value = -45367
$ws.Range("H120").Value = 21375
$ws.Range("J120").Value = 21375
$ws.Range("L120").Value = 21375
$ws.Range("N120").Value = -28633

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 23198.309
$ws.Range("I4").Value = 25089.834
$ws.Range("K4").Value = 75269.502
$ws.Range("M4").Value = -75157.502
$ws.Range("H6").Value = 12875.125
$ws.Range("J6").Value = 51000
$ws.Range("L6").Value = 153000
$ws.Range("N6").Value = -153226
$ws.Range("H54").Value = 3966.6667
$ws.Range("J54").Value = 3966.6667
$ws.Range("L54").Value = 11900.0001
$ws.Range("N54").Value = -13018.0001
$ws.Range("H92").Value = 3496.6
$ws.Range("I92").Value = 280
$ws.Range("J92").Value = 4300.75
$ws.Range("K92").Value = 840
$ws.Range("L92").Value = 12902.25
$ws.Range("M92").Value = 408
$ws.Range("N92").Value = -15398.25

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H9").Value = 52933
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 52933
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 52933
$ws.Range("M9").ClearContents()
$ws.Range("N9").Value = -53273
$ws.Range("H80").Value = 4138.4116
$ws.Range("I80").Value = 4250.8887
$ws.Range("J80").Value = 4011.875
$ws.Range("K80").Value = 4250.8887
$ws.Range("L80").Value = 4011.875
$ws.Range("M80").Value = -3252.8887
$ws.Range("N80").Value = -6007.875
$ws.Range("H83").Value = 4138.4116
$ws.Range("I83").Value = 4250.8887
$ws.Range("J83").Value = 4011.875
$ws.Range("K83").Value = 21254.4435
$ws.Range("L83").Value = 20059.375
$ws.Range("M83").Value = -16262.4435
$ws.Range("N83").Value = -30043.375
$ws.Range("H137").Value = 29587.273
$ws.Range("J137").Value = 29587.273
$ws.Range("L137").Value = 29587.273
$ws.Range("N137").Value = -39787.273

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2690.9092
$ws.Range("I7").Value = 1942.8572
$ws.Range("J7").Value = 4000
$ws.Range("K7").Value = 1942.8572
$ws.Range("L7").Value = 4000
$ws.Range("M7").Value = -1830.8572
$ws.Range("N7").Value = -4224
$ws.Range("H51").Value = 25063
$ws.Range("J51").Value = 25063
$ws.Range("L51").Value = 25063
$ws.Range("N51").Value = -26019
$ws.Range("H102").Value = 40000
$ws.Range("J102").Value = 40000
$ws.Range("L102").Value = 40000
$ws.Range("N102").Value = -46490
$ws.Range("H126").Value = 2690.9092
$ws.Range("I126").Value = 1942.8572
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 5828.571599999999
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -3358.571599999999
$ws.Range("N126").Value = -16940

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H102").Value = 40000
$ws.Range("J102").Value = 40000
$ws.Range("L102").Value = 40000
$ws.Range("N102").Value = -46490
$ws.Range("H138").Value = 29545.158
$ws.Range("J138").Value = 29545.158
$ws.Range("L138").Value = 29545.158
$ws.Range("N138").Value = -39825.158
